# Updated cryptos list — applies the price (column D) and volume(1h) (column E) changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.434.35"
$ws.Range("E2").Value = "  -3.00%  "

$ws.Range("D3").Value = "'3.184.42"
$ws.Range("E3").Value = "  -7.90%  "

$ws.Range("D5").Value = "'566.93"
$ws.Range("E5").Value = "  -3.69%  "

$ws.Range("D6").Value = "'170.80"
$ws.Range("E6").Value = "  -3.33%  "

$ws.Range("D7").Value = "'0.611"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'3.182.25"
$ws.Range("E9").Value = "  -7.90%  "

$ws.Range("E10").Value = "  -6.20%  "

$ws.Range("D11").Value = "'6.65"
$ws.Range("E11").Value = "  -4.49%  "

$ws.Range("D12").Value = "'0.397"
$ws.Range("E12").Value = "  -4.87%  "

$ws.Range("D13").Value = "'3.726.06"
$ws.Range("E13").Value = "  -8.16%  "

$ws.Range("E14").Value = "  +1.43%  "

$ws.Range("D15").Value = "'27.49"
$ws.Range("E15").Value = "  -9.38%  "

$ws.Range("D16").Value = "'64.398.64"
$ws.Range("E16").Value = "  -2.85%  "

$ws.Range("E17").Value = "  -5.09%  "

$ws.Range("D18").Value = "'3.182.38"
$ws.Range("E18").Value = "  -7.92%  "

$ws.Range("D19").Value = "'5.76"
$ws.Range("E19").Value = "  -3.48%  "

$ws.Range("D20").Value = "'13.03"
$ws.Range("E20").Value = "  -5.62%  "

$ws.Range("D21").Value = "'354.10"
$ws.Range("E21").Value = "  -5.28%  "

$ws.Range("D22").Value = "'7.21"
$ws.Range("E22").Value = "  -5.35%  "

$ws.Range("D24").Value = "'69.17"
$ws.Range("E24").Value = "  -5.70%  "

$ws.Range("D25").Value = "'0.0000121"
$ws.Range("E25").Value = "  -3.95%  "

$ws.Range("D26").Value = "'0.505"
$ws.Range("E26").Value = "  -5.59%  "

$ws.Range("D27").Value = "'9.57"
$ws.Range("E27").Value = "  -3.30%  "

$ws.Range("E28").Value = "  -0.91%  "

$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("E30").Value = "  -4.44%  "

$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("E32").Value = "  -4.96%  "

$ws.Range("D33").Value = "'22.11"
$ws.Range("E33").Value = "  -6.80%  "

$ws.Range("E34").Value = "  -5.07%  "

$ws.Range("D35").Value = "'6.67"
$ws.Range("E35").Value = "  -5.32%  "

$ws.Range("D36").Value = "'1.45"
$ws.Range("E36").Value = "  -6.48%  "

$ws.Range("E37").Value = "  -3.40%  "

$ws.Range("D38").Value = "'0.820"
$ws.Range("E38").Value = "  -7.49%  "

$ws.Range("E39").Value = "  -8.84%  "

$ws.Range("D40").Value = "'2.57"
$ws.Range("E40").Value = "  -1.80%  "

$ws.Range("D41").Value = "'1.71"
$ws.Range("E41").Value = "  -5.85%  "

$ws.Range("D42").Value = "'2.626.75"
$ws.Range("E42").Value = "  -4.99%  "

$ws.Range("D43").Value = "'4.19"
$ws.Range("E43").Value = "  -7.08%  "

$ws.Range("E44").Value = "  -6.27%  "

$ws.Range("D45").Value = "'39.66"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("D46").Value = "'0.0659"
$ws.Range("E46").Value = "  -5.04%  "

$ws.Range("D47").Value = "'23.91"
$ws.Range("E47").Value = "  -5.46%  "

$ws.Range("D48").Value = "'324.84"
$ws.Range("E48").Value = "  -4.29%  "

$ws.Range("E49").Value = "  -7.09%  "

$ws.Range("E50").Value = "  -0.85%  "

$ws.Range("D51").Value = "'0.999"
